$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("deep_linking")

$values = @(
    "/#/search?query=science&offsetIndex=0&searchType=ALL",
    "/#/search?query=science&offsetIndex=0&sortType=_score:desc&searchType=ARTICLES",
    "/#/search?query=science&offsetIndex=0&sortType=_score:desc&searchType=PATENTS",
    "/#/search?query=science&offsetIndex=0&sortType=_score:desc&searchType=PEOPLE",
    "/#/search?query=science&offsetIndex=0&sortType=sortdate:desc&searchType=POSTS",
    "/#/search?query=biology&offsetIndex=0&sortType=citingsrcslocalcount:desc&searchType=ALL",
    "/#/search?query=biology&offsetIndex=0&sortType=sortdate:desc&searchType=ALL",
    "/#/search?query=biology&offsetIndex=0&sortType=sortdate:asc&searchType=ALL",
    "/#/search?query=biology&offsetIndex=0&sortType=citingsrcslocalcount:desc&searchType=ARTICLES",
    "/#/search?query=biology&offsetIndex=0&sortType=sortdate:desc&searchType=ARTICLES",
    "/#/search?query=biology&offsetIndex=0&sortType=sortdate:asc&searchType=ARTICLES",
    "/#/search?query=biology&offsetIndex=0&sortType=citingsrcscount:desc&searchType=PATENTS",
    "/#/search?query=biology&offsetIndex=0&sortType=sortdate:desc&searchType=PATENTS",
    "/#/search?query=biology&offsetIndex=0&sortType=sortdate:asc&searchType=PATENTS",
    "/#/search?query=biology&offsetIndex=0&sortType=loadtime:desc&searchType=PEOPLE",
    "/#/search?query=post&offsetIndex=0&sortType=sortdate:desc&searchType=POSTS",
    "/#/search?query=biology&offsetIndex=0&sortType=sortdate:asc&searchType=POSTS",
    "/#/search?query=biology&offsetIndex=0&sortType=_score:desc&searchType=POSTS"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $values[$i]
}

$ws.Activate()
$ws.Range("A4").Select()
